# "Generate Report for handoff"
#
# The file's .md source is ready to hand off; this run produced the actual
# handoff package (a .xlf per target language) and the report needs to show
# that instead of the previous failure state:
#   - Status flips from "Handoff transform failed" to "Ready for handoff"
#     (Overview sheet + each per-language sheet all share that status text).
#   - Each per-language sheet gets its freshly produced "Latest Handoff File"
#     hyperlink + "Latest Handoff Datetime", and the "Handoff Reason" flips
#     from "Ignored" to "Include" now that the file is being sent out.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c838922efb9915e98f1ca7402c2b50c0687e2840"
$srcFile = "786c73be-f044-4ba5-8f37-ca032f9a9e06.md"
$readyStatus = "Ready for handoff"

# ---- Overview sheet: both language-status columns reflect the new status ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $readyStatus
$overview.Range("C2").Value = $readyStatus

# ---- Per-language detail sheets ----
$languages = @(
    @{ Sheet = "zh-cn"; Xlf = "786c73be-f044-4ba5-8f37-ca032f9a9e06.a4fd3d12f659c446c38962e04032f4fb75fc0d22.zh-cn.xlf"; Handoff = "2016-01-28 11:25:23" },
    @{ Sheet = "de-de"; Xlf = "786c73be-f044-4ba5-8f37-ca032f9a9e06.a4fd3d12f659c446c38962e04032f4fb75fc0d22.de-de.xlf"; Handoff = "2016-01-28 11:25:35" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status -> Ready for handoff
    $ws.Range("B2").Value = $readyStatus

    # Latest Handoff File -> new hyperlink to the generated .xlf, styled like
    # the workbook's other hyperlink cells (A2/A3 use the "HyperLink" style).
    $ws.Hyperlinks.Add($ws.Range("C2"), "$baseUrl/e2e/$($lang.Xlf)", "", "", $lang.Xlf)
    $ws.Range("C2").Style = "HyperLink"

    # Latest Handoff Datetime -> timestamp of this handoff run
    $ws.Range("D2").Value = $lang.Handoff

    # Handoff Reason -> Include (file is now part of the handoff), was Ignored
    $ws.Range("H2").Value = "Include"
}
